$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for data rows 2-43
# from 45847 (2025-07-09) to 45849 (2025-07-11), keeping existing formatting.
$ws.Range("C2:C43").Value = 45849
